$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/change (E) columns per latest cryptos snapshot
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.390.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.212.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.78%  "

$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "107.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -12.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "294.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.96%  "

$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.71%  "

$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.955"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.546.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.228.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.316.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +18.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "227.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.27%  "

$ws.Range("E32").Value = "  -5.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "173.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0882"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  -3.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0361"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.102"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.53%  "

$ws.Range("E42").Value = "  -5.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.229"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.20%  "

$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.56%  "

$ws.Range("E47").Value = "  -6.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.38%  "
